$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the missing U102 / TLC59711 row to the "Micro" BOM list (appended
# right after the existing last row of that column, row 23 -> row 24).
$ws.Range("A24").Value = "U102"
$ws.Range("B24").Value = "TLC59711"

# Move the view: scroll back to column A (no more topLeftCell override)
# and leave the selection on the newly added row.
$ws.Range("C24").Select() | Out-Null
